$p = $ppt.ActivePresentation

$s1 = $p.Slides.Add(1, 2)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = 'Learning with Experts for Fine-grained Category Discovery'
$s1.Shapes.Item(2).TextFrame.TextRange.Text = 'The paper "XCon: Learning with Experts for Fine-grained Category Discovery" by Yixin Fei et al. introduces a novel approach to address the problem of Generalized Category Discovery (GCD) in fine-grained datasets. GCD aims to cluster unlabeled images by leveraging information from a set of seen (labeled) classes, which is challenging due to large inter-class similarities and intra-class variances. XCon enhances the model''s ability to discern fine-grained discriminative features by partitioning the dataset into sub-datasets using k-means clustering and applying contrastive learning within these partitions.'

$s2 = $p.Slides.Add(2, 2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = 'Novel Category Discovery vs Fine-grained Category Discovery'
$s2.Shapes.Item(2).TextFrame.TextRange.Text = 'Earlier works in Novel Category Discovery (NCD) have utilized transfer learning and self-supervision techniques to categorize unseen classes based on knowledge transferred from seen classes. However, these methods are limited in their effectiveness in fine-grained classification tasks where class distinctions are more nuanced. XCon addresses this challenge by leveraging expert sub-datasets and a combination of contrastive learning approaches.'

$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = 'Partitioning the Dataset using k-means Clustering'
$s3.Shapes.Item(2).TextFrame.TextRange.Text = 'XCon partitions the dataset into several sub-datasets (referred to as expert datasets) using k-means clustering on self-supervised representations. This partitioning is based on the insight that class-irrelevant cues can lead to misleading clustering, thus, by creating more homogenized sub-datasets, the model is compelled to learn finer distinctions between classes.'

$s4 = $p.Slides.Add(4, 2)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = 'Supervised and Unsupervised Contrastive Learning'
$s4.Shapes.Item(2).TextFrame.TextRange.Text = 'The learning process involves both supervised and unsupervised contrastive learning across the full dataset and within the sub-datasets. The losses from these learning processes are combined to form the model''s final optimization objective. This dual approach allows XCon to learn from both coarse-grained (overall dataset) and fine-grained (sub-dataset) features.'

$s5 = $p.Slides.Add(5, 2)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = 'State-of-the-Art Performance on Fine-grained Category Discovery Benchmarks'
$s5.Shapes.Item(2).TextFrame.TextRange.Text = 'The effectiveness of XCon was evaluated on various datasets, including CIFAR-10/100, ImageNet-100, CUB-200, Stanford Cars, FGVC-Aircraft, and Oxford-IIIT Pet. These datasets encompass a wide range of classification challenges, from generic to fine-grained categories. XCon achieved state-of-the-art performance on several fine-grained category discovery benchmarks, demonstrating its capability to enhance classification performance significantly in challenging datasets.'

$s6 = $p.Slides.Add(6, 2)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = 'Validating the Effectiveness of XCon Components'
$s6.Shapes.Item(2).TextFrame.TextRange.Text = 'Ablation studies were conducted to validate the effectiveness of various components of XCon, including the impact of the weight of fine-grained loss and the number of sub-datasets. These studies confirmed the robustness of the XCon approach across different configurations and provided insights into optimal settings for different datasets.'

$s7 = $p.Slides.Add(7, 2)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = 'A Novel Approach to Fine-grained Category Discovery'
$s7.Shapes.Item(2).TextFrame.TextRange.Text = 'XCon introduces a novel method for fine-grained category discovery by leveraging expert sub-datasets and a combination of contrastive learning approaches. Its success across various benchmarks signifies a significant advancement in the field, offering a new direction for research in unsupervised and semi-supervised learning within fine-grained classification scenarios. The method''s codebase has been made publicly available, encouraging further exploration and adaptation of this approach within the research community.'
